$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last refreshed" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 01:09"

# 2) Update numeric stats for the countries whose figures changed.
#    (Column layout: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#     E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes)

# Estados Unidos
$ws.Range("B4").Value = 5793942
$ws.Range("C4").Value = 47670
$ws.Range("D4").Value = 3117184
$ws.Range("E4").Value = 2497658
$ws.Range("G4").Value = 1044
$ws.Range("H4").Value = 179100

# Brasil
$ws.Range("D5").Value = 2670755
$ws.Range("E5").Value = 748217

# Alemania
$ws.Range("B23").Value = 233021
$ws.Range("C23").Value = 1737
$ws.Range("E23").Value = 17893

# Canada
$ws.Range("B27").Value = 124282
$ws.Range("C27").Value = 409
$ws.Range("D27").Value = 110604
$ws.Range("E27").Value = 4616
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 9062

# Egipto
$ws.Range("B34").Value = 97148
$ws.Range("C34").Value = 123
$ws.Range("D34").Value = 64318
$ws.Range("E34").Value = 27599
$ws.Range("G34").Value = 19
$ws.Range("H34").Value = 5231

# Japon
$ws.Range("B48").Value = 59721
$ws.Range("C48").Value = 1220
$ws.Range("D48").Value = 46467
$ws.Range("E48").Value = 12099
$ws.Range("G48").Value = 11
$ws.Range("H48").Value = 1155

# Nigeria
$ws.Range("B52").Value = 51304
$ws.Range("C52").Value = 340
$ws.Range("D52").Value = 37885
$ws.Range("E52").Value = 12423
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 996

# Kirguistan
$ws.Range("E56").Value = 5818
$ws.Range("H56").Value = 1054

# Venezuela overtakes Uzbekistan, Afganistan and Etiopia in the ranking,
# so rows 60-63 shift: Venezuela's updated numbers move into row 60, and
# Uzbekistan / Afganistan / Etiopia each drop one row (their own figures
# are unchanged, only their position moves).
$ws.Range("A60").Value = "Venezuela"
$ws.Range("B60").Value = 38219
$ws.Range("C60").Value = 652
$ws.Range("D60").Value = 27306
$ws.Range("E60").Value = 10596
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 317

$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 38074
$ws.Range("C61").Value = 527
$ws.Range("D61").Value = 33989
$ws.Range("E61").Value = 3825
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 260

$ws.Range("A62").Value = "Afganistan"
$ws.Range("B62").Value = 37894
$ws.Range("C62").Value = 38
$ws.Range("D62").Value = 28016
$ws.Range("E62").Value = 8493
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 1385

$ws.Range("A63").Value = "Etiopia"
$ws.Range("B63").Value = 37665
$ws.Range("C63").Value = 1829
$ws.Range("D63").Value = 13913
$ws.Range("E63").Value = 23115
$ws.Range("G63").Value = 17
$ws.Range("H63").Value = 637

# Chequia
$ws.Range("B74").Value = 21551
$ws.Range("C74").Value = 506
$ws.Range("E74").Value = 5100
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 411

# Noruega
$ws.Range("B89").Value = 10275
$ws.Range("C89").Value = 78
$ws.Range("E89").Value = 861

# Guayana Francesa
$ws.Range("B94").Value = 8777
$ws.Range("C94").Value = 34
$ws.Range("D94").Value = 8251
$ws.Range("E94").Value = 471
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 55

# Gabon
$ws.Range("B95").Value = 8388
$ws.Range("C95").Value = 69
$ws.Range("D95").Value = 6734
$ws.Range("E95").Value = 1601

# Luxemburgo
$ws.Range("B101").Value = 7704
$ws.Range("C101").Value = 67
$ws.Range("D101").Value = 6969
$ws.Range("E101").Value = 611

# Guinea Ecuatorial
$ws.Range("B110").Value = 4926
$ws.Range("C110").Value = 34
$ws.Range("E110").Value = 2130

# Guyana
$ws.Range("B163").Value = 881
$ws.Range("C163").Value = 35
$ws.Range("D163").Value = 433
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 30

# Trinidad yTobago
$ws.Range("B164").Value = 864
$ws.Range("C164").Value = 97
$ws.Range("E164").Value = 690
